$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '54.504.52'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +5.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.173.51'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.14%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '399.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '110.37'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.28%  '
$ws.Range("E7").Value = '  +0.55%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  +3.78%  '
$ws.Range("E10").Value = '  +4.99%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0883'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.679.71'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.08'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.61%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.02'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.69%  '
$ws.Range("E16").Value = '  +7.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.186.30'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.80%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.55'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '54.415.48'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.81%  '
$ws.Range("E20").Value = '  +2.53%  '
$ws.Range("B21").Value = 'InternetComputer(DFINITY)'
$ws.Range("C21").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.89'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.93%  '
$ws.Range("B22").Value = 'ShibaInu'
$ws.Range("C22").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0000100'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.44%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.14'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '275.07'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.86%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.97'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '27.78'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.57'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.09%  '
$ws.Range("E29").Value = '  -0.54%  '
$ws.Range("E30").Value = '  -0.11%  '
$ws.Range("E31").Value = '  +2.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.02'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.27%  '
$ws.Range("E33").Value = '  +13.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '36.55'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.76%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.10'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '51.55'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.36%  '
$ws.Range("E37").Value = '  +5.98%  '
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.89'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +10.55%  '
$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.09'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +10.04%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.94'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.292'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.20'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.63%  '
$ws.Range("E44").Value = '  +1.69%  '
$ws.Range("E45").Value = '  +0.97%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.05'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.48'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.82%  '
$ws.Range("E48").Value = '  -0.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.091.46'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.85%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0341'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0510'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +14.21%  '
